$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark additional sprint progress (hours logged today) ---
# H9: 1 hour logged against T01x task row -> highlight like other logged cells
$ws.Range("H9").Value = 1
$ws.Range("H9").Interior.Color = 65535

# L12 / L14: 0.5 hour logged -> highlight like other logged cells in that column
$ws.Range("L12").Value = 0.5
$ws.Range("L12").Interior.Color = 65535

$ws.Range("L14").Value = 0.5
$ws.Range("L14").Interior.Color = 65535

# --- Burndown total for the final day now accounts for the scope change ---
# (difference between the originally estimated hours and the hours actually tracked)
$ws.Range("V32").Formula = "=U32-SUM(V3:V29)-(SUM(B3:B29)-SUM(C3:V29))"

# --- Update the saved selection / scroll position ---
[void]$ws.Range("Z30").Select()
